$d = $word.ActiveDocument

$replacements = @(
    @{old = "588÷9=65, 3"; new = "604÷4=151, 0"},
    @{old = "189÷9=21, 0"; new = "814÷2=407, 0"},
    @{old = "317÷5=63, 2"; new = "904÷9=100, 4"},
    @{old = "979÷6=163, 1"; new = "930÷8=116, 2"},
    @{old = "312÷2=156, 0"; new = "714÷9=79, 3"},
    @{old = "602÷7=86, 0"; new = "256÷5=51, 1"},
    @{old = "846÷4=211, 2"; new = "583÷6=97, 1"},
    @{old = "131÷4=32, 3"; new = "745÷9=82, 7"},
    @{old = "633÷4=158, 1"; new = "187÷8=23, 3"},
    @{old = "331÷9=36, 7"; new = "105÷8=13, 1"},
    @{old = "556÷4=139, 0"; new = "983÷4=245, 3"},
    @{old = "332÷3=110, 2"; new = "503÷4=125, 3"},
    @{old = "466÷6=77, 4"; new = "536÷2=268, 0"},
    @{old = "764÷7=109, 1"; new = "856÷6=142, 4"},
    @{old = "558÷2=279, 0"; new = "516÷9=57, 3"},
    @{old = "325÷6=54, 1"; new = "470÷5=94, 0"},
    @{old = "321÷8=40, 1"; new = "879÷7=125, 4"},
    @{old = "695÷3=231, 2"; new = "710÷5=142, 0"},
    @{old = "991÷7=141, 4"; new = "555÷2=277, 1"},
    @{old = "728÷8=91, 0"; new = "243÷8=30, 3"},
    @{old = "639÷6=106, 3"; new = "535÷6=89, 1"},
    @{old = "858÷4=214, 2"; new = "502÷6=83, 4"},
    @{old = "260÷2=130, 0"; new = "505÷2=252, 1"},
    @{old = "281÷8=35, 1"; new = "512÷9=56, 8"},
    @{old = "156÷7=22, 2"; new = "843÷4=210, 3"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
